# Updated cryptos list data: prices and 1h volume deltas, plus two
# row re-orderings (rows 13/14 and 42/43 swap their Coin/Link/Price/Volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.892.36'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  +0.48%  '
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = '''1.639.41'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +0.83%  '
$ws.Range('E3').Style = 'Normal'

$ws.Range('E4').Value = '''  +0.22%  '
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = '''215.81'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.55%  '
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = '''0.5078'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.18%  '
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Value = '''1.004'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.28%  '
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = '''0.2599'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  +1.60%  '
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Value = '''0.06460'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +1.75%  '
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = '''20.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +4.91%  '
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = '''0.07840'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.95%  '
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Value = '''4.267'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +0.27%  '
$ws.Range('E12').Style = 'Normal'

$ws.Range('B13').Value = '''WrappedEther'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = '''1.642.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +1.04%  '
$ws.Range('E13').Style = 'Normal'

$ws.Range('B14').Value = '''WrappedliquidstakedEther2.0'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = '''1.863.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  +0.78%  '
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = '''0.5655'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.94%  '
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = '''0.0₅7704'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +2.86%  '
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = '''63.45'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -0.29%  '
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = '''25.888.81'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +0.40%  '
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = '''1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.25%  '
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = '''194.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.48%  '
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = '''4.386'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -0.46%  '
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = '''9.988'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  +2.30%  '
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = '''6.199'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +3.76%  '
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = '''1.004'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +0.27%  '
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Value = '''1.770'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -5.21%  '
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = '''138.17'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -1.96%  '
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Value = '''0.1236'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  -0.11%  '
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = '''6.867'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  +2.01%  '
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = '''15.61'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  +1.22%  '
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = '''1.246'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.76%  '
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = '''0.05015'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  +2.75%  '
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = '''3.312'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  +0.58%  '
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = '''3.252'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +2.33%  '
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = '''1.576'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +1.88%  '
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').Value = '''  +1.06%  '
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Value = '''0.9075'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +1.53%  '
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').Value = '''2.582'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +1.87%  '
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').Value = '''1.133.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +0.06%  '
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Value = '''0.5526'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +0.42%  '
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = '''0.01575'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +1.37%  '
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Value = '''0.9943'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.65%  '
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = '''Quant'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''99.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +2.70%  '
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').Value = '''FraxShare'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''5.499'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.05%  '
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = '''0.8013'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +0.98%  '
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').Value = '''  -2.19%  '
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = '''55.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +1.86%  '
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = '''0.4242'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -3.99%  '
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').Value = '''7.687'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  +1.07%  '
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').Value = '''0.05042'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -1.54%  '
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = '''1.001'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  +0.08%  '
$ws.Range('E50').Style = 'Normal'

$ws.Range('E51').Value = '''  +0.12%  '
$ws.Range('E51').Style = 'Normal'
